$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.184.12"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "2.517.52"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "109.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "321.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0819"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "2.912.25"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "2.516.18"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.852"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "48.036.04"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.46%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E30").Value = "  +6.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "2.035.14"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.71%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.38%  "
